# Append 9 new rows (regcntr_id/machine_id pairs) to the
# master-reg_center_machine test-data sheet, then leave the sheet in the
# same "just finished typing at the bottom" view state the original author
# left it in (selection on the row below the new data, scrolled down),
# and turn on portrait page setup, matching the committed workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# regcntr_id / machine_id pairs for the new rows; lang_code, is_active,
# cr_by and cr_dtimes repeat the same values used by every existing row.
$newRows = @(
    @(10002, 10021),
    @(10003, 10022),
    @(10004, 10023),
    @(10005, 10024),
    @(10006, 10025),
    @(10007, 10026),
    @(10008, 10027),
    @(10009, 10028),
    @(10010, 10029)
)

$startRow = 22
$r = $startRow
foreach ($pair in $newRows) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $r = $r + 1
}

# Match the author's final on-screen state: everything below the newly
# typed data selected (as Excel does when you select the remaining blank
# rows after entering data), with the view scrolled near the bottom of
# the table.
$lastDataRow = $r - 1
$selRowsRange = ($lastDataRow + 1).ToString() + ":1048576"
$ws.Rows($selRowsRange).Select()

# Page setup was switched to portrait orientation in the committed file.
$ws.PageSetup.Orientation = 1

Write-Host "Added $($newRows.Count) rows; last data row = $lastDataRow"
